$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May")

# Make sure the May sheet (already the active/tab-selected one) is active
$ws.Activate()

# Fill in hours-paid ($50) for the first three logged entries
$ws.Range("D2").Value = 50
$ws.Range("D3").Value = 50
$ws.Range("D4").Value = 50

# Add the new row 5 entry: task, hours, fee
$ws.Range("B5").Value = "Decoding transcription, editing Jane's story, sending out"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 150

# Row 5 now wraps onto two lines like the other long task descriptions
$ws.Rows.Item(5).RowHeight = 36

# Move the selection to match where the user ended up
$ws.Range("F5").Select() | Out-Null
